$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: new backlog item (plain style, E column gets the "tip" highlight fill
# used on rows 10/12/14) ---
$ws.Range("A16").Value = 15
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "Bewoners"
$ws.Range("E16").Value = "Ik wil een veiligheidsknop die op elk moment te bereiken is die linkt naar een pagina met de meest voorkomende remedies tegen ongelukken. 1e hulp, wat te doen bij brand etc. Dit wil ik om de veiligheid te bevorderen."

$ws.Range("E10").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 17: follow-up backlog item, shaded like the alternating rows
# (3/5/7/9/11/13/15) ---
$ws.Range("A17").Value = 16
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "Bewoners"
$ws.Range("E17").Value = "Ik wil op elke stap van het koken een nuttige tip hebben die de veiligheid en / of efficientie van het koken bevorderd, denk hierbij aan natte doek onder de snijplank om weg glijden te voorkomen etc. Dit wil ik om de veiligheid te bevorderen."

$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Column E widened slightly to fit the new content
$ws.Columns.Item(5).ColumnWidth = 213.5

# Selection moves on, matching the saved view state
$ws.Range("E28").Select() | Out-Null
